# Weekly fruit/vegetable data refresh: two new price observations were
# collected (new rows), pushing the existing rows down by two.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows right above the first data row (row 17),
# shifting every existing record down by two rows.
$ws.Rows.Item(17).Insert()
$ws.Rows.Item(17).Insert()

# The "constant" columns are identical for every record in this sheet.
$marketId = 7
$market   = "Terminal Hortofrutícola Agro Chillán"
$region   = "Ñuble"
$codreg   = 16
$catId    = 100112037
$category = "Cebollín"
$variety  = "Sin especificar"
$classification = "Hortaliza"

# New row 17
$ws.Range("A17").Value = $marketId
$ws.Range("B17").Value = $market
$ws.Range("C17").Value = $region
$ws.Range("D17").Value = 44883
$ws.Range("E17").Value = $codreg
$ws.Range("F17").Value = $catId
$ws.Range("G17").Value = $category
$ws.Range("H17").Value = $variety
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 400
$ws.Range("K17").Value = 600
$ws.Range("L17").Value = 700
$ws.Range("M17").Value = 650
$ws.Range("N17").Value = "$/paquete 6 unidades"
$ws.Range("O17").Value = "Provincia de Diguillín"
$ws.Range("P17").Value = 108
$ws.Range("Q17").Value = 6
$ws.Range("R17").Value = $classification

# New row 18
$ws.Range("A18").Value = $marketId
$ws.Range("B18").Value = $market
$ws.Range("C18").Value = $region
$ws.Range("D18").Value = 44883
$ws.Range("E18").Value = $codreg
$ws.Range("F18").Value = $catId
$ws.Range("G18").Value = $category
$ws.Range("H18").Value = $variety
$ws.Range("I18").Value = "Segunda"
$ws.Range("J18").Value = 300
$ws.Range("K18").Value = 500
$ws.Range("L18").Value = 500
$ws.Range("M18").Value = 500
$ws.Range("N18").Value = "$/paquete 6 unidades"
$ws.Range("O18").Value = "Provincia de Diguillín"
$ws.Range("P18").Value = 83
$ws.Range("Q18").Value = 6
$ws.Range("R18").Value = $classification

# Keep the date columns formatted the same way as every other row
# (the row-insert already copies formatting, but make sure explicitly).
$ws.Range("D17").NumberFormat = $ws.Range("D19").NumberFormat
$ws.Range("D18").NumberFormat = $ws.Range("D19").NumberFormat
